# Update cryptocurrency price/volume data as of Wed Feb  8 07:44:15 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'332.23"
$ws.Range("E2").Value = "'1.10%"
$ws.Range("D3").Value = "'45.88"
$ws.Range("E3").Value = "'4.29%"
$ws.Range("D4").Value = "'5.605"
$ws.Range("E4").Value = "'1.58%"
$ws.Range("E5").Value = "'4.32%"
$ws.Range("D6").Value = "'2.041"
$ws.Range("E6").Value = "'2.42%"
$ws.Range("D7").Value = "'0.9744"
$ws.Range("E7").Value = "'2.58%"
$ws.Range("E8").Value = "'-0.55%"
$ws.Range("D9").Value = "'0.1157"
$ws.Range("E9").Value = "'2.85%"
$ws.Range("D10").Value = "'0.1915"
$ws.Range("E10").Value = "'1.40%"
$ws.Range("D11").Value = "'10.36"
$ws.Range("E11").Value = "'-3.17%"
$ws.Range("D12").Value = "'0.09922"
$ws.Range("E12").Value = "'-0.60%"
$ws.Range("D13").Value = "'0.04691"
$ws.Range("E13").Value = "'-2.29%"
$ws.Range("D14").Value = "'0.1060"
$ws.Range("E14").Value = "'-0.26%"
$ws.Range("D15").Value = "'0.001290"
$ws.Range("E15").Value = "'1.40%"
$ws.Range("D16").Value = "'0.006068"
$ws.Range("E16").Value = "'1.57%"
$ws.Range("E17").Value = "'0.38%"
$ws.Range("D18").Value = "'4.453"
$ws.Range("E18").Value = "'1.79%"
$ws.Range("D19").Value = "'0.3366"
$ws.Range("E19").Value = "'-3.12%"
$ws.Range("E20").Value = "'-1.92%"
$ws.Range("D21").Value = "'0.2652"
$ws.Range("E21").Value = "'4.15%"
$ws.Range("D22").Value = "'0.04195"
$ws.Range("E22").Value = "'3.02%"
$ws.Range("D23").Value = "'0.001312"
$ws.Range("E23").Value = "'3.50%"
$ws.Range("D24").Value = "'0.004579"
$ws.Range("E24").Value = "'5.90%"
$ws.Range("E25").Value = "'8.46%"
$ws.Range("D26").Value = "'0.0003747"
$ws.Range("E26").Value = "'0.01%"
$ws.Range("D38").Value = "'0.02758"
$ws.Range("E38").Value = "'6.55%"
$ws.Range("D39").Value = "'0.05788"
$ws.Range("E39").Value = "'2.39%"
$ws.Range("D40").Value = "'0.007687"
$ws.Range("E40").Value = "'1.76%"
$ws.Range("D41").Value = "'0.1435"
$ws.Range("E41").Value = "'2.82%"
$ws.Range("D42").Value = "'0.007295"
$ws.Range("E42").Value = "'-1.47%"
$ws.Range("D43").Value = "'0.002014"
$ws.Range("E43").Value = "'-0.14%"
$ws.Range("D44").Value = "'0.008195"
$ws.Range("E44").Value = "'-5.17%"
$ws.Range("D45").Value = "'0.3404"
$ws.Range("D46").Value = "'0.00007292"
$ws.Range("E46").Value = "'2.47%"
$ws.Range("E47").Value = "'0.12%"
$ws.Range("D48").Value = "'0.0005812"
$ws.Range("E48").Value = "'0.01%"
$ws.Range("E49").Value = "'-7.35%"
$ws.Range("D50").Value = "'0.003505"
$ws.Range("E50").Value = "'-0.78%"
$ws.Range("D51").Value = "'0.00002104"
$ws.Range("E51").Value = "'0.12%"
